# Update cryptocurrency price (D) and 1h volume change (E) columns
# to match the refreshed data feed, per the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 2: D2: '26.652.95' -> '26.669.13'; E2: '  +1.44%  ' -> '  +1.43%  '
$ws.Cells.Item(2, 4).Value = "26.669.13"
$ws.Cells.Item(2, 5).Value = "  +1.43%  "

# Row 3: D3: '1.633.77' -> '1.634.82'; E3: '  +0.95%  ' -> '  +0.96%  '
$ws.Cells.Item(3, 4).Value = "1.634.82"
$ws.Cells.Item(3, 5).Value = "  +0.96%  "

# Row 4: E4: '  +0.10%  ' -> '  +0.06%  '
$ws.Cells.Item(4, 5).Value = "  +0.06%  "

# Row 5: D5: '213.25' -> '213.37'; E5: '  +0.63%  ' -> '  +0.68%  '
$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = "213.37"
$ws.Cells.Item(5, 4).Style = "Normal"
$ws.Cells.Item(5, 5).Value = "  +0.68%  "

# Row 6: D6: '0.503' -> '0.502'; E6: '  +3.81%  ' -> '  +3.83%  '
$ws.Cells.Item(6, 4).NumberFormat = "@"
$ws.Cells.Item(6, 4).Value = "0.502"
$ws.Cells.Item(6, 4).Style = "Normal"
$ws.Cells.Item(6, 5).Value = "  +3.83%  "

# Row 7: E7: '  +0.09%  ' -> '  +0.08%  '
$ws.Cells.Item(7, 5).Value = "  +0.08%  "

# Row 8: E8: '  +2.48%  ' -> '  +2.39%  '
$ws.Cells.Item(8, 5).Value = "  +2.39%  "

# Row 9: E9: '  +1.45%  ' -> '  +1.46%  '
$ws.Cells.Item(9, 5).Value = "  +1.46%  "

# Row 10: E10: '  +2.80%  ' -> '  +2.40%  '
$ws.Cells.Item(10, 5).Value = "  +2.40%  "

# Row 11: E11: '  +3.49%  ' -> '  +3.55%  '
$ws.Cells.Item(11, 5).Value = "  +3.55%  "

# Row 12: D12: '1.862.63' -> '1.861.78'; E12: '  +1.07%  ' -> '  +0.96%  '
$ws.Cells.Item(12, 4).Value = "1.861.78"
$ws.Cells.Item(12, 5).Value = "  +0.96%  "

# Row 13: D13: '1.666.13' -> '1.650.26'; E13: '  +2.97%  ' -> '  +2.55%  '
$ws.Cells.Item(13, 4).Value = "1.650.26"
$ws.Cells.Item(13, 5).Value = "  +2.55%  "

# Row 14: E14: '  +2.59%  ' -> '  +2.40%  '
$ws.Cells.Item(14, 5).Value = "  +2.40%  "

# Row 15: D15: '0.527' -> '0.526'; E15: '  +1.86%  ' -> '  +1.63%  '
$ws.Cells.Item(15, 4).NumberFormat = "@"
$ws.Cells.Item(15, 4).Value = "0.526"
$ws.Cells.Item(15, 4).Style = "Normal"
$ws.Cells.Item(15, 5).Value = "  +1.63%  "

# Row 16: D16: '26.644.29' -> '26.660.02'; E16: '  +1.36%  ' -> '  +1.37%  '
$ws.Cells.Item(16, 4).Value = "26.660.02"
$ws.Cells.Item(16, 5).Value = "  +1.37%  "

# Row 17: D17: '63.49' -> '63.43'; E17: '  +2.06%  ' -> '  +1.90%  '
$ws.Cells.Item(17, 4).NumberFormat = "@"
$ws.Cells.Item(17, 4).Value = "63.43"
$ws.Cells.Item(17, 4).Style = "Normal"
$ws.Cells.Item(17, 5).Value = "  +1.90%  "

# Row 18: D18: '0.0₃0744' -> '0.0₃0743'; E18: '  +2.51%  ' -> '  +2.26%  '
$ws.Cells.Item(18, 4).Value = "0.0₃0743"
$ws.Cells.Item(18, 5).Value = "  +2.26%  "

# Row 19: D19: '218.80' -> '219.76'; E19: '  +8.73%  ' -> '  +9.06%  '
$ws.Cells.Item(19, 4).NumberFormat = "@"
$ws.Cells.Item(19, 4).Value = "219.76"
$ws.Cells.Item(19, 4).Style = "Normal"
$ws.Cells.Item(19, 5).Value = "  +9.06%  "

# Row 20: E20: '  +0.02%  ' -> '  +0.08%  '
$ws.Cells.Item(20, 5).Value = "  +0.08%  "

# Row 21: E21: '  +0.59%  ' -> '  +0.48%  '
$ws.Cells.Item(21, 5).Value = "  +0.48%  "

# Row 22: D22: '9.46' -> '9.47'; E22: '  +1.66%  ' -> '  +1.56%  '
$ws.Cells.Item(22, 4).NumberFormat = "@"
$ws.Cells.Item(22, 4).Value = "9.47"
$ws.Cells.Item(22, 4).Style = "Normal"
$ws.Cells.Item(22, 5).Value = "  +1.56%  "

# Row 23: E23: '  +2.91%  ' -> '  +2.80%  '
$ws.Cells.Item(23, 5).Value = "  +2.80%  "

# Row 24: E24: '  +2.10%  ' -> '  +1.85%  '
$ws.Cells.Item(24, 5).Value = "  +1.85%  "

# Row 25: D25: '148.52' -> '148.77'; E25: '  +2.84%  ' -> '  +2.95%  '
$ws.Cells.Item(25, 4).NumberFormat = "@"
$ws.Cells.Item(25, 4).Value = "148.77"
$ws.Cells.Item(25, 4).Style = "Normal"
$ws.Cells.Item(25, 5).Value = "  +2.95%  "

# Row 26: E26: '  +0.11%  ' -> '  +0.02%  '
$ws.Cells.Item(26, 5).Value = "  +0.02%  "

# Row 27: E27: '  +1.36%  ' -> '  +1.65%  '
$ws.Cells.Item(27, 5).Value = "  +1.65%  "

# Row 28: D28: '6.92' -> '6.93'; E28: '  +5.91%  ' -> '  +5.72%  '
$ws.Cells.Item(28, 4).NumberFormat = "@"
$ws.Cells.Item(28, 4).Value = "6.93"
$ws.Cells.Item(28, 4).Style = "Normal"
$ws.Cells.Item(28, 5).Value = "  +5.72%  "

# Row 29: D29: '15.50' -> '15.60'; E29: '  +2.36%  ' -> '  +2.98%  '
$ws.Cells.Item(29, 4).NumberFormat = "@"
$ws.Cells.Item(29, 4).Value = "15.60"
$ws.Cells.Item(29, 4).Style = "Normal"
$ws.Cells.Item(29, 5).Value = "  +2.98%  "

# Row 30: E30: '  -1.58%  ' -> '  -1.88%  '
$ws.Cells.Item(30, 5).Value = "  -1.88%  "

# Row 31: E31: '  -0.23%  ' -> '  +0.12%  '
$ws.Cells.Item(31, 5).Value = "  +0.12%  "

# Row 32: E32: '  +4.38%  ' -> '  +4.31%  '
$ws.Cells.Item(32, 5).Value = "  +4.31%  "

# Row 33: E33: '  +2.90%  ' -> '  +2.69%  '
$ws.Cells.Item(33, 5).Value = "  +2.69%  "

# Row 34: E34: '  +1.45%  ' -> '  +1.19%  '
$ws.Cells.Item(34, 5).Value = "  +1.19%  "

# Row 35: E35: '  -0.04%  ' -> '  -0.03%  '
$ws.Cells.Item(35, 5).Value = "  -0.03%  "

# Row 36: D36: '1.208.36' -> '1.213.91'; E36: '  +2.74%  ' -> '  +3.06%  '
$ws.Cells.Item(36, 4).Value = "1.213.91"
$ws.Cells.Item(36, 5).Value = "  +3.06%  "

# Row 37: D37: '0.0172' -> '0.0173'; E37: '  +5.38%  ' -> '  +5.58%  '
$ws.Cells.Item(37, 4).NumberFormat = "@"
$ws.Cells.Item(37, 4).Value = "0.0173"
$ws.Cells.Item(37, 4).Style = "Normal"
$ws.Cells.Item(37, 5).Value = "  +5.58%  "

# Row 38: D38: '0.811' -> '0.813'; E38: '  +0.87%  ' -> '  +0.88%  '
$ws.Cells.Item(38, 4).NumberFormat = "@"
$ws.Cells.Item(38, 4).Value = "0.813"
$ws.Cells.Item(38, 4).Style = "Normal"
$ws.Cells.Item(38, 5).Value = "  +0.88%  "

# Row 40: D40: '0.505' -> '0.506'; E40: '  +2.14%  ' -> '  +2.11%  '
$ws.Cells.Item(40, 4).NumberFormat = "@"
$ws.Cells.Item(40, 4).Value = "0.506"
$ws.Cells.Item(40, 4).Style = "Normal"
$ws.Cells.Item(40, 5).Value = "  +2.11%  "

# Row 41: E41: '  -1.08%  ' -> '  -1.16%  '
$ws.Cells.Item(41, 5).Value = "  -1.16%  "

# Row 42: D42: '5.44' -> '5.43'; E42: '  +1.87%  ' -> '  +1.75%  '
$ws.Cells.Item(42, 4).NumberFormat = "@"
$ws.Cells.Item(42, 4).Value = "5.43"
$ws.Cells.Item(42, 4).Style = "Normal"
$ws.Cells.Item(42, 5).Value = "  +1.75%  "

# Row 43: D43: '0.792' -> '0.795'; E43: '  +0.48%  ' -> '  +0.76%  '
$ws.Cells.Item(43, 4).NumberFormat = "@"
$ws.Cells.Item(43, 4).Value = "0.795"
$ws.Cells.Item(43, 4).Style = "Normal"
$ws.Cells.Item(43, 5).Value = "  +0.76%  "

# Row 44: D44: '1.770.96' -> '1.770.59'; E44: '  +0.96%  ' -> '  +0.88%  '
$ws.Cells.Item(44, 4).Value = "1.770.59"
$ws.Cells.Item(44, 5).Value = "  +0.88%  "

# Row 45: D45: '93.26' -> '93.46'; E45: '  +0.84%  ' -> '  +0.82%  '
$ws.Cells.Item(45, 4).NumberFormat = "@"
$ws.Cells.Item(45, 4).Value = "93.46"
$ws.Cells.Item(45, 4).Style = "Normal"
$ws.Cells.Item(45, 5).Value = "  +0.82%  "

# Row 46: E46: '  +1.97%  ' -> '  +1.51%  '
$ws.Cells.Item(46, 5).Value = "  +1.51%  "

# Row 47: D47: '54.81' -> '54.86'; E47: '  +2.26%  ' -> '  +2.08%  '
$ws.Cells.Item(47, 4).NumberFormat = "@"
$ws.Cells.Item(47, 4).Value = "54.86"
$ws.Cells.Item(47, 4).Style = "Normal"
$ws.Cells.Item(47, 5).Value = "  +2.08%  "

# Row 48: E48: '  +1.06%  ' -> '  +1.02%  '
$ws.Cells.Item(48, 5).Value = "  +1.02%  "

# Row 49: D49: '7.73' -> '7.69'; E49: '  +6.28%  ' -> '  +5.87%  '
$ws.Cells.Item(49, 4).NumberFormat = "@"
$ws.Cells.Item(49, 4).Value = "7.69"
$ws.Cells.Item(49, 4).Style = "Normal"
$ws.Cells.Item(49, 5).Value = "  +5.87%  "

# Row 50: E50: '  +0.51%  ' -> '  +0.37%  '
$ws.Cells.Item(50, 5).Value = "  +0.37%  "

# Row 51: E51: '  +0.30%  ' -> '  +0.22%  '
$ws.Cells.Item(51, 5).Value = "  +0.22%  "
